$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("ФИО" etc. shift right by one)
$ws.Columns("C").Insert()

# New header for the inserted column
$ws.Range("C1").Value = "InviteID"

# New invite id value (numeric) for row 2
$ws.Range("C2").Value = 351029552

# Replace the old "Кем_работает"/"Компания" data with new invite-related
# numbers, kept as text (they are long digit strings, same as the
# existing numberStoredAsText convention used in this sheet)
$ws.Range("E2").Value = "'11111111111"
$ws.Range("F2").Value = "'22222222222"
$ws.Range("G2").Value = "'111111111111"
